$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format (one cell at a time -- a comma-separated union range
# only applies NumberFormat to its first area in this host) on price cells
# whose new values would otherwise be auto-coerced to numbers by Excel.
$textCells = @('D5', 'D6', 'D9', 'D10', 'D12', 'D13', 'D14', 'D17', 'D19', 'D21', 'D22', 'D23', 'D26', 'D28', 'D29', 'D30', 'D32', 'D34', 'D35', 'D36', 'D38', 'D40', 'D44', 'D45', 'D47', 'D48', 'D50')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '42.897.90'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '2.296.16'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '300.52'
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = '97.02'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.506'
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('D10').Value = '35.76'
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = '17.86'
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '0.117'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '6.77'
$ws.Range('E14').Value = '  -1.92%  '
$ws.Range('D15').Value = '2.649.44'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '2.291.70'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').Value = '0.776'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '42.810.79'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').Value = '12.86'
$ws.Range('E19').Value = '  -4.54%  '
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D21').Value = '6.05'
$ws.Range('E21').Value = '  -2.15%  '
$ws.Range('D22').Value = '67.77'
$ws.Range('E22').Value = '  -0.30%  '
$ws.Range('D23').Value = '240.83'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').Value = '  -0.76%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '25.21'
$ws.Range('E28').Value = '  +1.13%  '
$ws.Range('D29').Value = '165.70'
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').Value = '2.03'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('D32').Value = '32.87'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').Value = '1.00'
$ws.Range('D35').Value = '5.02'
$ws.Range('E35').Value = '  -4.19%  '
$ws.Range('D36').Value = '17.00'
$ws.Range('E36').Value = '  -7.54%  '
$ws.Range('E37').Value = '  -1.62%  '
$ws.Range('D38').Value = '0.0686'
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('D40').Value = '1.76'
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('E42').Value = '  -1.67%  '
$ws.Range('D43').Value = '2.018.91'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').Value = '0.0282'
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('D45').Value = '10.11'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('D47').Value = '17.20'
$ws.Range('E47').Value = '  -1.80%  '
$ws.Range('D48').Value = '2.79'
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('D50').Value = '53.34'
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('D51').Value = '2.514.54'
$ws.Range('E51').Value = '  -0.94%  '
